# Update the "想去人数" (want-to-go count) column F values on both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same rows).
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 266
    5  = 837
    6  = 6
    7  = 293
    8  = 7476
    11 = 127
    12 = 102
    15 = 17
    18 = 248
    19 = 677
    20 = 18
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
